$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("hsplogin")
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws1.Range("B3").Value = "k1577"

[void]$ws1.Activate()
[void]$ws1.Range("B3").Select()

[void]$ws2.Activate()
[void]$ws2.Range("L7").Select()
